$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: was 228.6118008523737 -> now blank
$ws.Range("D3").ClearContents()

# C4: 15857.7055216777 -> 20.76103376777028
$ws.Range("C4").Value = 20.76103376777028

# C5: 64729.72874398364 -> 0
$ws.Range("C5").Value = 0

# Row 7 label changes from "Other" to "Biogas", and D7 gets a value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 100.5891923750444

# New row 8: "Other" row, same style as row 7's A cell, D8 = 0
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = "Other"
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 0
